$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying report rows (2-8) have been re-sorted/re-fetched; most of the
# per-row "fixed" attributes (location, coordinates, dates, observers, ...)
# stay the same, but the species identity (columns A,B,D,E,F,G,H) and the
# start/end time (columns Z,AB) move to different rows. One record additionally
# got a new report id (A) and updated times (Z/AB) because it was re-submitted.

$rows = @(
  @{ Row = 2; A = 111708888; B = 90678; D = "LC"; E = 4366; F = "Skarp dropptaggsvamp";   G = "Hydnellum peckii";      H = "Banker";                  Z = "13:54"; AB = "13:54" },
  @{ Row = 3; A = 111708920; B = 90666; D = "LC"; E = 4364; F = "Dropptaggsvamp";          G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst.";     Z = "13:53"; AB = "13:53" },
  @{ Row = 4; A = 111704319; B = 90710; D = "NT"; E = 5449; F = "Svart taggsvamp";         G = "Phellodon niger";      H = "(Fr.:Fr.) P.Karst.";      Z = "15:11"; AB = "15:11" },
  @{ Row = 5; A = 111706580; B = 88032; D = "VU"; E = 6276; F = "Goliatmusseron";          G = "Tricholoma matsutake"; H = "(S.Ito & S.Imai) Singer"; Z = "14:48"; AB = "14:48" },
  @{ Row = 6; A = 111708040; B = 90658; D = "NT"; E = 4361; F = "Orange taggsvamp";        G = "Hydnellum aurantiacum"; H = "(Batsch:Fr.) P.Karst.";  Z = "14:20"; AB = "14:20" },
  @{ Row = 7; A = 111708099; B = 90660; D = "NT"; E = 4362; F = "Blå taggsvamp";           G = "Hydnellum caeruleum";  H = "(Hornem.) P.Karst.";      Z = "14:16"; AB = "14:16" },
  @{ Row = 8; A = 111708029; B = 90662; D = "LC"; E = 4363; F = "Zontaggsvamp";            G = "Hydnellum concrescens"; H = "(Pers.) Banker";         Z = "14:21"; AB = "14:21" }
)

foreach ($r in $rows) {
  $rowNum = $r.Row
  $ws.Cells.Item($rowNum, 1).Value = $r.A    # A - Id
  $ws.Cells.Item($rowNum, 2).Value = $r.B    # B - Taxonsorteringsordning
  $ws.Cells.Item($rowNum, 4).Value = $r.D    # D - Rödlistade
  $ws.Cells.Item($rowNum, 5).Value = $r.E    # E - TaxonId
  $ws.Cells.Item($rowNum, 6).Value = $r.F    # F - Artnamn
  $ws.Cells.Item($rowNum, 7).Value = $r.G    # G - Vetenskapligt namn
  $ws.Cells.Item($rowNum, 8).Value = $r.H    # H - Auktor
  $ws.Cells.Item($rowNum, 26).Value = $r.Z   # Z - Starttid
  $ws.Cells.Item($rowNum, 28).Value = $r.AB  # AB - Sluttid
}
